$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (pushes existing rows 6..62 down to 7..63,
# mirroring the dimension growing from A1:R62 to A1:R63).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new daily price record.
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44699
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100112032
$ws.Range("G6").Value = "Zapallo italiano"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("N6").Value = "$/caja 60 unidades"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 208
$ws.Range("Q6").Value = 60
$ws.Range("R6").Value = "Hortaliza"
